$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.513.54'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.485.17'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.54'
$ws.Range("E5").Value = '  -2.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.83'
$ws.Range("E6").Value = '  -3.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.157'
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.42'
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").Value = '2.923.99'
$ws.Range("D14").Value = '58.391.42'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.45'
$ws.Range("E15").Value = '  -3.71%  '
$ws.Range("D17").Value = '2.477.60'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.94'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.22'
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '322.20'
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  -1.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.36'
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("E26").Value = '  -3.20%  '
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").Value = '0.0₃0753'
$ws.Range("E28").Value = '  -3.55%  '
$ws.Range("E29").Value = '  -4.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.71'
$ws.Range("E30").Value = '  -4.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.92'
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -4.73%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.25'
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.35'
$ws.Range("E36").Value = '  -8.22%  '
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.798'
$ws.Range("E39").Value = '  -3.72%  '
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '278.10'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.98'
$ws.Range("E42").Value = '  -5.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.596'
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '127.29'
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0914'
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0497'
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '17.26'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").Value = '1.743.99'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.972'
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.66'
$ws.Range("E51").Value = '  -2.18%  '
